$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.093771920020204
$ws.Range("D2").Value = 1.103523098914172
$ws.Range("E2").Value = 1.09671933217123
$ws.Range("F2").Value = 1.110981167979094
$ws.Range("I2").Value = 1.049359121094061
$ws.Range("J2").Value = 1.098586239806329
$ws.Range("K2").Value = 1.106137551252725
$ws.Range("L2").Value = 1.099350889984954
$ws.Range("M2").Value = 1.113577128896065
$ws.Range("N2").Value = 1.100146358276456

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.09560557136484
$ws.Range("D3").Value = 1.105301676400755
$ws.Range("E3").Value = 1.098388855592123
$ws.Range("F3").Value = 1.11281134670521
$ws.Range("I3").Value = 1.049739421189455
$ws.Range("J3").Value = 1.100082184469998
$ws.Range("K3").Value = 1.107736396463868
$ws.Range("L3").Value = 1.100839680383622
$ws.Range("M3").Value = 1.115228814501649
$ws.Range("N3").Value = 1.101644427353136

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.096789243614594
$ws.Range("D4").Value = 1.10645002574066
$ws.Range("E4").Value = 1.099466703382474
$ws.Range("F4").Value = 1.113993239271897
$ws.Range("I4").Value = 1.049982936447796
$ws.Range("J4").Value = 1.101046983498589
$ws.Range("K4").Value = 1.108767980415643
$ws.Range("L4").Value = 1.101800096800784
$ws.Range("M4").Value = 1.116294749700577
$ws.Range("N4").Value = 1.102610596507011

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.097286199363747
$ws.Range("D5").Value = 1.106932205569278
$ws.Range("E5").Value = 1.099919258624457
$ws.Range("F5").Value = 1.114489557772938
$ws.Range("I5").Value = 1.050084700294026
$ws.Range("J5").Value = 1.101451837286973
$ws.Range("K5").Value = 1.109200959408902
$ws.Range("L5").Value = 1.102203166919671
$ws.Range("M5").Value = 1.116742208526509
$ws.Range("N5").Value = 1.103016025234212

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.09736960216833
$ws.Range("D6").Value = 1.107013131725524
$ws.Range("E6").Value = 1.09999521149011
$ws.Range("F6").Value = 1.114572860058335
$ws.Range("I6").Value = 1.050101751236428
$ws.Range("J6").Value = 1.101519770545437
$ws.Range("K6").Value = 1.109273617893696
$ws.Range("L6").Value = 1.102270804131172
$ws.Range("M6").Value = 1.116817300508385
$ws.Range("N6").Value = 1.103084054965695

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.096795886533728
$ws.Range("D7").Value = 1.106456470934775
$ws.Range("E7").Value = 1.099472752676711
$ws.Range("F7").Value = 1.113999873241479
$ws.Range("I7").Value = 1.049984298612604
$ws.Range("J7").Value = 1.101052396095065
$ws.Range("K7").Value = 1.108773768629807
$ws.Range("L7").Value = 1.10180548533319
$ws.Range("M7").Value = 1.116300731243058
$ws.Range("N7").Value = 1.102616016789995

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.094392205051339
$ws.Range("D8").Value = 1.104124705806202
$ws.Range("E8").Value = 1.097284069938292
$ws.Range("F8").Value = 1.111600182084229
$ws.Range("I8").Value = 1.049488178636169
$ws.Range("J8").Value = 1.099092466970454
$ws.Range("K8").Value = 1.106678512808549
$ws.Range("L8").Value = 1.099854648089623
$ws.Range("M8").Value = 1.114135915148425
$ws.Range("N8").Value = 1.10065330434122

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.090134283863673
$ws.Range("D9").Value = 1.099995967472548
$ws.Range("E9").Value = 1.093407993699381
$ws.Range("F9").Value = 1.107352914299869
$ws.Range("I9").Value = 1.048594134221807
$ws.Range("J9").Value = 1.095613906090884
$ws.Range("K9").Value = 1.102963003940098
$ws.Range("L9").Value = 1.096394009922563
$ws.Range("M9").Value = 1.110299045348919
$ws.Range("N9").Value = 1.097169803506205

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.087279663656874
$ws.Range("D10").Value = 1.09722919880705
$ws.Range("E10").Value = 1.090810087071864
$ws.Range("F10").Value = 1.104507891268124
$ws.Range("I10").Value = 1.047984529558414
$ws.Range("J10").Value = 1.093277296202349
$ws.Range("K10").Value = 1.100469408827325
$ws.Range("L10").Value = 1.094070644560377
$ws.Range("M10").Value = 1.107725341703821
$ws.Range("N10").Value = 1.094829875363615

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.086039562675418
$ws.Range("D11").Value = 1.096027566023009
$ws.Range("E11").Value = 1.089681687140609
$ws.Range("F11").Value = 1.103272551004686
$ws.Range("I11").Value = 1.04771728776628
$ws.Range("J11").Value = 1.092261168003457
$ws.Range("K11").Value = 1.099385528185156
$ws.Range("L11").Value = 1.093060562723907
$ws.Range("M11").Value = 1.106606958876841
$ws.Range("N11").Value = 1.093812304146129

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.08557830898608
$ws.Range("D12").Value = 1.095580666907205
$ws.Range("E12").Value = 1.089262008369981
$ws.Range("F12").Value = 1.102813157708313
$ws.Range("I12").Value = 1.047617524705941
$ws.Range("J12").Value = 1.091883061716865
$ws.Range("K12").Value = 1.098982288419987
$ws.Range("L12").Value = 1.092684749452861
$ws.Range("M12").Value = 1.106190930987332
$ws.Range("N12").Value = 1.093433660905241

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.08567727806106
$ws.Range("D13").Value = 1.095676553877099
$ws.Range("E13").Value = 1.089352055631184
$ws.Range("F13").Value = 1.102911723538719
$ws.Range("I13").Value = 1.047638946816537
$ws.Range("J13").Value = 1.09196419742137
$ws.Range("K13").Value = 1.099068813868815
$ws.Range("L13").Value = 1.092765391159572
$ws.Range("M13").Value = 1.106280198291952
$ws.Range("N13").Value = 1.093514911831753

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.086001448163006
$ws.Range("D14").Value = 1.095990636713743
$ws.Range("E14").Value = 1.089647007491138
$ws.Range("F14").Value = 1.103234588388474
$ws.Range("I14").Value = 1.047709051499426
$ws.Range("J14").Value = 1.092229927401227
$ws.Range("K14").Value = 1.099352209377019
$ws.Range("L14").Value = 1.09302951069967
$ws.Range("M14").Value = 1.106572582438762
$ws.Range("N14").Value = 1.093781019178661

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.086201096667081
$ws.Range("D15").Value = 1.096184078880895
$ws.Range("E15").Value = 1.089828664906976
$ws.Range("F15").Value = 1.103433444915599
$ws.Range("I15").Value = 1.047752179228125
$ws.Range("J15").Value = 1.092393563015594
$ws.Range("K15").Value = 1.099526733629983
$ws.Range("L15").Value = 1.093192160324141
$ws.Range("M15").Value = 1.106752648526062
$ws.Range("N15").Value = 1.09394488717437

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.087361878494166
$ws.Range("D16").Value = 1.097308869778544
$ws.Range("E16").Value = 1.090884900334464
$ws.Range("F16").Value = 1.104589802997924
$ws.Range("I16").Value = 1.048002196012501
$ws.Range("J16").Value = 1.093344640057809
$ws.Range("K16").Value = 1.100541253829041
$ws.Range("L16").Value = 1.094137593728649
$ws.Range("M16").Value = 1.107799480379393
$ws.Range("N16").Value = 1.094897314855073

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.088088913360699
$ws.Range("D17").Value = 1.098013444090497
$ws.Range("E17").Value = 1.091546503105616
$ws.Range("F17").Value = 1.105314225766705
$ws.Range("I17").Value = 1.048158143603267
$ws.Range("J17").Value = 1.093940046919861
$ws.Range("K17").Value = 1.101176516144986
$ws.Range("L17").Value = 1.094729544096421
$ws.Range("M17").Value = 1.108455060336619
$ws.Range("N17").Value = 1.095493567263164

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.088512592819958
$ws.Range("D18").Value = 1.098424063873489
$ws.Range("E18").Value = 1.091932069279172
$ws.Range("F18").Value = 1.105736439936995
$ws.Range("I18").Value = 1.048248789159975
$ws.Range("J18").Value = 1.094286918007205
$ws.Range("K18").Value = 1.101546655829505
$ws.Range("L18").Value = 1.095074429090797
$ws.Range("M18").Value = 1.108837068861887
$ws.Range("N18").Value = 1.09584093094724

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.088656991389469
$ws.Range("D19").Value = 1.098564016371872
$ws.Range("E19").Value = 1.092063480967502
$ws.Range("F19").Value = 1.105880348735639
$ws.Range("I19").Value = 1.048279643495691
$ws.Range("J19").Value = 1.09440512137991
$ws.Range("K19").Value = 1.101672796875942
$ws.Range("L19").Value = 1.095191960338989
$ws.Range("M19").Value = 1.108967259922876
$ws.Range("N19").Value = 1.095959302182292

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.088010949638509
$ws.Range("D20").Value = 1.097937885902566
$ws.Range("E20").Value = 1.091475554214663
$ws.Range("F20").Value = 1.105236536259042
$ws.Range("I20").Value = 1.048141444627208
$ws.Range("J20").Value = 1.093876208920672
$ws.Range("K20").Value = 1.101108399805443
$ws.Range("L20").Value = 1.094666073838598
$ws.Range("M20").Value = 1.108384762256844
$ws.Range("N20").Value = 1.095429638606694

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.085906005555755
$ws.Range("D21").Value = 1.095898162727067
$ws.Range("E21").Value = 1.08956016658252
$ws.Range("F21").Value = 1.103139527595172
$ws.Range("I21").Value = 1.047688421192647
$ws.Range("J21").Value = 1.092151695197514
$ws.Range("K21").Value = 1.099268774216149
$ws.Range("L21").Value = 1.092951751426598
$ws.Range("M21").Value = 1.10648649959068
$ws.Range("N21").Value = 1.093702675876245

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.084578918082301
$ws.Range("D22").Value = 1.094612463511199
$ws.Range("E22").Value = 1.088352749137761
$ws.Range("F22").Value = 1.101817962584928
$ws.Range("I22").Value = 1.047400706162709
$ws.Range("J22").Value = 1.091063534596851
$ws.Range("K22").Value = 1.098108427590563
$ws.Range("L22").Value = 1.09187027134551
$ws.Range("M22").Value = 1.105289444587298
$ws.Range("N22").Value = 1.092612969962715

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.085282782414362
$ws.Range("D23").Value = 1.095294350503609
$ws.Range("E23").Value = 1.088993126926664
$ws.Range("F23").Value = 1.102518848433047
$ws.Range("I23").Value = 1.047553504113009
$ws.Range("J23").Value = 1.091640763233616
$ws.Range("K23").Value = 1.09872390570254
$ws.Range("L23").Value = 1.092443932534641
$ws.Range("M23").Value = 1.105924367534182
$ws.Range("N23").Value = 1.093191018330353

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.088046179286267
$ws.Range("D24").Value = 1.097972028465745
$ws.Range("E24").Value = 1.091507614001093
$ws.Range("F24").Value = 1.105271641820969
$ws.Range("I24").Value = 1.048148991152172
$ws.Range("J24").Value = 1.093905055861148
$ws.Range("K24").Value = 1.101139179877064
$ws.Range("L24").Value = 1.094694754519826
$ws.Range("M24").Value = 1.108416528107639
$ws.Range("N24").Value = 1.095458526513134

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.091237806185844
$ws.Range("D25").Value = 1.10106579380668
$ws.Range("E25").Value = 1.094412431307
$ws.Range("F25").Value = 1.108453247130737
$ws.Range("I25").Value = 1.048827640258023
$ws.Range("J25").Value = 1.096516233684216
$ws.Range("K25").Value = 1.103926412677228
$ws.Range("L25").Value = 1.097291476923027
$ws.Range("M25").Value = 1.11129368676199
$ws.Range("N25").Value = 1.098073412508217
